$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: Remove the "Meta description" paragraph entirely.
# It sits right after the title (Heading1) paragraph and reads:
#   "Meta description" (bold) + ": Fowl Play London offers a London and
#   Victorian era theme, engaging mini-game, and improved graphics.
#   Read our review and play for free."
# ---------------------------------------------------------------------
$metaRange = $d.Content
$metaFound = $metaRange.Find.Execute(
    "Meta description: Fowl Play London offers a London and Victorian era theme, engaging mini-game, and improved graphics. Read our review and play for free.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($metaFound) {
    $metaRange.Expand(4) | Out-Null  # wdParagraph - grab the whole paragraph incl. mark
    $metaRange.Delete() | Out-Null
}

# ---------------------------------------------------------------------
# Step 2: Insert a new bold paragraph ("Play Fowl Play London Free: A
# Victorian London Slot Game Review") right before the final paragraph
# that currently holds the italic "Prompt: ..." image-prompt text.
# ---------------------------------------------------------------------
$promptRange = $d.Content
$promptFound = $promptRange.Find.Execute(
    "Prompt: Create a feature image for Fowl Play London",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($promptFound) {
    $promptRange.Expand(4) | Out-Null  # wdParagraph
    $insertPos = $promptRange.Start

    $insertPoint = $d.Range($insertPos, $insertPos)
    $insertPoint.InsertParagraphBefore() | Out-Null

    $newParaRange = $d.Range($insertPos, $insertPos)
    $newParaRange.Expand(4) | Out-Null  # wdParagraph - the freshly created empty paragraph

    $newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fowl Play London Free: A Victorian London Slot Game Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $newParaRange.InsertXML($newParaXml) | Out-Null
}

# ---------------------------------------------------------------------
# Step 3: Replace the text of the (still italic) final paragraph with
# the repurposed meta-description copy.
# ---------------------------------------------------------------------
$oldPromptText = "Prompt: Create a feature image for Fowl Play London, a cartoon-style image featuring a happy Maya warrior with glasses. The image should include the game's title and showcase the London theme. The image should have a fun and engaging vibe that reflects the game's updated gameplay. Use bold and vibrant colors to catch the reader's eye and make the image stand out. The Maya warrior should be dressed in a Sherlock Holmes outfit, with a magnifying glass in hand, investigating the London scenery in the background. Make sure the image is high-resolution and in a landscape format suitable for online articles."
$newPromptText = "Fowl Play London offers a London and Victorian era theme, engaging mini-game, and improved graphics. Read our review and play for free."

$d.Content.Find.Execute($oldPromptText, $true, $false, $false, $false, $false,
                         $true, 1, $false, $newPromptText, 2) | Out-Null
